# Fix Training Data Issue (#48)
# Data was taken from 1 day off due to way NBA stats were shown.
# The "Date" column (BF) contained the mangled literal text "5-12-2011-12"
# for every data row; correct it to the proper ISO-style date text
# "2012-05-12". The values must remain plain text (not get converted to a
# real Excel date serial number), so we force text with a leading
# apostrophe.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 31; $row++) {
    $cell = $ws.Cells.Item($row, 58)  # column BF = 58
    if ($cell.Value() -eq "5-12-2011-12") {
        # Leading apostrophe forces the text to stay a literal string
        # instead of being auto-recognized/converted to a date serial.
        $cell.Value = "'2012-05-12"
    }
}
